$d = $word.ActiveDocument

$d.Content.Find.Execute("Reflection 1: {{Q1Text}}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Reflection 1: {Q1Text}", 2)
$d.Content.Find.Execute("Reflection 2: {{Q2Text}}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Reflection 2: {Q2Text}", 2)
$d.Content.Find.Execute("Reflection 3: {{Q3Text}}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Reflection 3: {Q3Text}", 2)
